# Bump the "Förändrad" (Changed) date in column C for every data row
# (rows 2-418) from 2026-02-21 (serial 46074) to 2026-02-22 (serial 46075).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C418").Value = 46075
